$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.87207273792896
$ws.Range("C2").Value = 8.16641005421905
$ws.Range("D2").Value = 5.984966522721663
$ws.Range("E2").Value = 10.78010055468211
$ws.Range("G2").Value = 46.32993351566965
$ws.Range("H2").Value = 17.37094115793137
$ws.Range("M2").Value = 15.3123776500617
$ws.Range("B3").Value = 13.36998827415124
$ws.Range("C3").Value = 7.638791975103773
$ws.Range("D3").Value = 5.867586134613891
$ws.Range("E3").Value = 10.64564611964875
$ws.Range("G3").Value = 45.10422381062273
$ws.Range("H3").Value = 17.24139484836551
$ws.Range("M3").Value = 15.07321361649098
$ws.Range("B4").Value = 13.05710524186207
$ws.Range("C4").Value = 7.294899043597051
$ws.Range("D4").Value = 5.796330801672381
$ws.Range("E4").Value = 10.56543734702103
$ws.Range("G4").Value = 44.34875625077297
$ws.Range("H4").Value = 17.16543885217095
$ws.Range("M4").Value = 14.92972147265454
$ws.Range("B5").Value = 12.92868484488207
$ws.Range("C5").Value = 7.149736715679578
$ws.Range("D5").Value = 5.767541948358864
$ws.Range("E5").Value = 10.53337420973583
$ws.Range("G5").Value = 44.04064078631554
$ws.Range("H5").Value = 17.13540720025811
$ws.Range("M5").Value = 14.87217035709392
$ws.Range("B6").Value = 12.9073126650564
$ws.Range("C6").Value = 7.125328947994826
$ws.Range("D6").Value = 5.762777830853788
$ws.Range("E6").Value = 10.52808867079981
$ws.Range("G6").Value = 43.98947694894016
$ws.Range("H6").Value = 17.13047662348134
$ws.Range("M6").Value = 14.86267207346279
$ws.Range("B7").Value = 13.05537670605684
$ws.Range("C7").Value = 7.292961685494484
$ws.Range("D7").Value = 5.795941483108913
$ws.Range("E7").Value = 10.56500236974079
$ws.Range("G7").Value = 44.34460130593416
$ws.Range("H7").Value = 17.16503008162546
$ws.Range("M7").Value = 14.92894147716251
$ws.Range("B8").Value = 13.70004609745212
$ws.Range("C8").Value = 7.988617555513635
$ws.Range("D8").Value = 5.944349952953877
$ws.Range("E8").Value = 10.73327352714496
$ws.Range("G8").Value = 45.90818054307304
$ws.Range("H8").Value = 17.32554145700036
$ws.Range("M8").Value = 15.22926552957215
$ws.Range("B9").Value = 14.91827127274751
$ws.Range("C9").Value = 9.195192595211925
$ws.Range("D9").Value = 6.239989789387902
$ws.Range("E9").Value = 11.08037957688356
$ws.Range("G9").Value = 48.93200225388976
$ws.Range("H9").Value = 17.66782661392969
$ws.Range("M9").Value = 15.84127496019544
$ws.Range("B10").Value = 15.77406403455126
$ws.Range("C10").Value = 9.986431391047878
$ws.Range("D10").Value = 6.457560793257739
$ws.Range("E10").Value = 11.34376807556432
$ws.Range("G10").Value = 51.10329379817811
$ws.Range("H10").Value = 17.93473053189931
$ws.Range("M10").Value = 16.3002042163727
$ws.Range("B11").Value = 16.15289990630977
$ws.Range("C11").Value = 10.32590323215413
$ws.Range("D11").Value = 6.556125652046261
$ws.Range("E11").Value = 11.46495538247506
$ws.Range("G11").Value = 52.07549217042104
$ws.Range("H11").Value = 18.0591854399492
$ws.Range("M11").Value = 16.50999789795102
$ws.Range("B12").Value = 16.2947114398872
$ws.Range("C12").Value = 10.45152509930231
$ws.Range("D12").Value = 6.593351978228343
$ws.Range("E12").Value = 11.51100484973151
$ws.Range("G12").Value = 52.44105516338664
$ws.Range("H12").Value = 18.10672187900767
$ws.Range("M12").Value = 16.58950734206634
$ws.Range("B13").Value = 16.26424501348027
$ws.Range("C13").Value = 10.42460020799405
$ws.Range("D13").Value = 6.585339585903117
$ws.Range("E13").Value = 11.50108083683405
$ws.Range("G13").Value = 52.362445032834
$ws.Range("H13").Value = 18.09646637718502
$ws.Range("M13").Value = 16.57238194534216
$ws.Range("B14").Value = 16.16460054294675
$ws.Range("C14").Value = 10.33629683209116
$ws.Range("D14").Value = 6.559190466582808
$ws.Range("E14").Value = 11.46874094410512
$ws.Range("G14").Value = 52.10562103691402
$ws.Range("H14").Value = 18.06308828193757
$ws.Range("M14").Value = 16.51653834417046
$ws.Range("B15").Value = 16.10334729064615
$ws.Range("C15").Value = 10.28182748555938
$ws.Range("D15").Value = 6.543159478718392
$ws.Range("E15").Value = 11.44895132164413
$ws.Range("G15").Value = 51.94796196033175
$ws.Range("H15").Value = 18.04269550354589
$ws.Range("M15").Value = 16.4823386117638
$ws.Range("B16").Value = 15.74908319458798
$ws.Range("C16").Value = 9.963835483583809
$ws.Range("D16").Value = 6.451107716252155
$ws.Range("E16").Value = 11.33587247190019
$ws.Range("G16").Value = 51.0394148756847
$ws.Range("H16").Value = 17.92665580996482
$ws.Range("M16").Value = 16.28650723917201
$ws.Range("B17").Value = 15.52896715029302
$ws.Range("C17").Value = 9.763525773934855
$ws.Range("D17").Value = 6.394503601714079
$ws.Range("E17").Value = 11.26682574297961
$ws.Range("G17").Value = 50.47781132215012
$ws.Range("H17").Value = 17.85622727581915
$ws.Range("M17").Value = 16.16657418963265
$ws.Range("B18").Value = 15.40138408143262
$ws.Range("C18").Value = 9.646385626283228
$ws.Range("D18").Value = 6.361910112686028
$ws.Range("E18").Value = 11.22724319882004
$ws.Range("G18").Value = 50.15334908057289
$ws.Range("H18").Value = 17.81600648411059
$ws.Range("M18").Value = 16.09769307931029
$ws.Range("B19").Value = 15.35802320523196
$ws.Range("C19").Value = 9.606392538481064
$ws.Range("D19").Value = 6.350869509628307
$ws.Range("E19").Value = 11.21386496249924
$ws.Range("G19").Value = 50.04325520704877
$ws.Range("H19").Value = 17.80243871831541
$ws.Range("M19").Value = 16.07439104524665
$ws.Range("B20").Value = 15.5525011828743
$ws.Range("C20").Value = 9.785048449983979
$ws.Range("D20").Value = 6.400533253470984
$ws.Range("E20").Value = 11.27416257587872
$ws.Range("G20").Value = 50.53774695167276
$ws.Range("H20").Value = 17.86369493012349
$ws.Range("M20").Value = 16.17933139610805
$ws.Range("B21").Value = 16.19391422465256
$ws.Range("C21").Value = 10.36231304378146
$ws.Range("D21").Value = 6.566874060042094
$ws.Range("E21").Value = 11.47823596010503
$ws.Range("G21").Value = 52.18112936750479
$ws.Range("H21").Value = 18.07288139764554
$ws.Range("M21").Value = 16.5329398382542
$ws.Range("B22").Value = 16.60347115179036
$ws.Range("C22").Value = 10.72252468082287
$ws.Range("D22").Value = 6.674999185803117
$ws.Range("E22").Value = 11.61251572770498
$ws.Range("G22").Value = 53.23994962609076
$ws.Range("H22").Value = 18.21196138453742
$ws.Range("M22").Value = 16.7643868013367
$ws.Range("B23").Value = 16.38580713686814
$ws.Range("C23").Value = 10.53182950774857
$ws.Range("D23").Value = 6.617357009069834
$ws.Range("E23").Value = 11.54077785403328
$ws.Range("G23").Value = 52.67633980505268
$ws.Range("H23").Value = 18.13752519408211
$ws.Range("M23").Value = 16.6408542509122
$ws.Range("B24").Value = 15.54186466104608
$ws.Range("C24").Value = 9.77532421469801
$ws.Range("D24").Value = 6.397807403861817
$ws.Range("E24").Value = 11.27084523776339
$ws.Range("G24").Value = 50.51065497540343
$ws.Range("H24").Value = 17.86031796196293
$ws.Range("M24").Value = 16.17356363877131
$ws.Range("B25").Value = 14.59486750423561
$ws.Range("C25").Value = 8.885569963318767
$ws.Range("D25").Value = 6.159774749083089
$ws.Range("E25").Value = 10.98485559958379
$ws.Range("G25").Value = 48.12114044747106
$ws.Range("H25").Value = 17.57241661556716
$ws.Range("M25").Value = 15.67373484681342
